$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws4 = $wb.Worksheets.Item(4)

# Sheet 1
$ws1.Range("F3").Value = 310
$ws1.Range("F4").Value = 2959
$ws1.Range("F5").Value = 76
$ws1.Range("F8").Value = 1650
$ws1.Range("F10").Value = 847
$ws1.Range("F11").Value = 111
$ws1.Range("F12").Value = 4
$ws1.Range("G12").Value = "不可售"
$ws1.Range("F13").Value = 2654
$ws1.Range("F15").Value = 1505
$ws1.Range("F16").Value = 7019
$ws1.Range("F18").Value = 7176
$ws1.Range("F20").Value = 5383
$ws1.Range("F21").Value = 3100
$ws1.Range("F22").Value = 3469
$ws1.Range("F23").Value = 223
$ws1.Range("F24").Value = 168
$ws1.Range("F25").Value = 1861
$ws1.Range("F26").Value = 79
$ws1.Range("F30").Value = 174
$ws1.Range("F32").Value = 2398
$ws1.Range("F33").Value = 1156
$ws1.Range("F34").Value = 2644
$ws1.Range("F35").Value = 18
$ws1.Range("F37").Value = 164
$ws1.Range("F38").Value = 380
$ws1.Range("F39").Value = 1051
$ws1.Range("F41").Value = 470

# Sheet 2
$ws2.Range("F14").Value = 94
$ws2.Range("F15").Value = 17
$ws2.Range("F17").Value = 60
$ws2.Range("F18").Value = 5

# Sheet 4
$ws4.Range("F4").Value = 310
$ws4.Range("F6").Value = 2959
$ws4.Range("F7").Value = 76
$ws4.Range("F9").Value = 1650
$ws4.Range("F11").Value = 847
$ws4.Range("F12").Value = 111
$ws4.Range("F14").Value = 2654
$ws4.Range("F15").Value = 1505
$ws4.Range("F19").Value = 7019
$ws4.Range("F21").Value = 7176
$ws4.Range("F23").Value = 5383
$ws4.Range("F24").Value = 3100
$ws4.Range("F25").Value = 3469
$ws4.Range("F27").Value = 224
$ws4.Range("F29").Value = 1861
$ws4.Range("F30").Value = 17
$ws4.Range("F35").Value = 174
$ws4.Range("F37").Value = 2398
$ws4.Range("F38").Value = 1156
$ws4.Range("F39").Value = 60
$ws4.Range("F40").Value = 2644
$ws4.Range("F41").Value = 18
$ws4.Range("F43").Value = 164
$ws4.Range("F44").Value = 5
$ws4.Range("F45").Value = 380
$ws4.Range("F46").Value = 1051
$ws4.Range("F48").Value = 470
